$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 55
$ws.Range("H55").Value = 2676.6667
$ws.Range("I55").Value = 974.5
$ws.Range("J55").Value = 3527.75
$ws.Range("K55").Value = 974.5
$ws.Range("L55").Value = 3527.75
$ws.Range("M55").Value = -760.5
$ws.Range("N55").Value = -3955.75

# Row 62
$ws.Range("H62").Value = 4989.3335
$ws.Range("J62").Value = 4989.3335
$ws.Range("L62").Value = 4989.3335
$ws.Range("N62").Value = -6237.3335

# Row 65
$ws.Range("H65").Value = 4989.3335
$ws.Range("J65").Value = 4989.3335
$ws.Range("L65").Value = 24946.6675
$ws.Range("N65").Value = -31186.6675

# Row 106
$ws.Range("H106").Value = 6388.6665
$ws.Range("I106").Value = 5197.4
$ws.Range("K106").Value = 5197.4
$ws.Range("M106").Value = -4566.4

# Row 135
$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

# Row 138
$ws.Range("H138").Value = 2818.3635
$ws.Range("J138").Value = 4475.4
$ws.Range("L138").Value = 13426.2
$ws.Range("N138").Value = -23706.2


# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 492.44446
$ws.Range("I2").Value = 390.57144
$ws.Range("K2").Value = 390.57144
$ws.Range("M2").Value = -277.57144

# Row 45
$ws.Range("H45").Value = 1901.5
$ws.Range("I45").Value = 1417.6
$ws.Range("K45").Value = 1417.6
$ws.Range("M45").Value = -1040.6

# Row 74
$ws.Range("H74").Value = 488
$ws.Range("I74").Value = 488
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 488
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = 386
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 488
$ws.Range("I77").Value = 488
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2440
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = 1928
$ws.Range("M77").ClearContents()

# Row 116
$ws.Range("H116").Value = 492.44446
$ws.Range("I116").Value = 390.57144
$ws.Range("K116").Value = 390.57144
$ws.Range("M116").Value = 1903.42856

# Row 122
$ws.Range("H122").Value = 2666.6667
$ws.Range("I122").Value = 1500
$ws.Range("K122").Value = 4500
$ws.Range("M122").Value = -2050

# Row 132
$ws.Range("H132").Value = 3017.5
$ws.Range("I132").Value = 1426.9
$ws.Range("K132").Value = 4280.700000000001
$ws.Range("M132").Value = -1750.700000000001


# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 492.44446
$ws.Range("I3").Value = 390.57144
$ws.Range("K3").Value = 390.57144
$ws.Range("M3").Value = -276.57144

# Row 94
$ws.Range("H94").Value = 689.2
$ws.Range("I94").Value = 724
$ws.Range("J94").Value = 550
$ws.Range("K94").Value = 724
$ws.Range("L94").Value = 550
$ws.Range("M94").Value = -273
$ws.Range("N94").Value = -1452

# Row 107
$ws.Range("H107").Value = 4013
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4013
$ws.Range("K107").Value = 0
$ws.Range("M107").Value = 4013
$ws.Range("N107").Value = -7853
$ws.Range("L107").ClearContents()

# Row 134
$ws.Range("H134").Value = 1832.4
$ws.Range("I134").Value = 850.55554
$ws.Range("K134").Value = 2551.66662
$ws.Range("M134").Value = -16.66661999999997


# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 2570.1853
$ws.Range("I31").Value = 1375.3529
$ws.Range("J31").Value = 4601.4
$ws.Range("K31").Value = 1375.3529
$ws.Range("L31").Value = 4601.4
$ws.Range("M31").Value = -1080.3529
$ws.Range("N31").Value = -5191.4

# Row 34
$ws.Range("H34").Value = 2570.1853
$ws.Range("I34").Value = 1375.3529
$ws.Range("J34").Value = 4601.4
$ws.Range("K34").Value = 1375.3529
$ws.Range("L34").Value = 4601.4
$ws.Range("M34").Value = -1173.3529
$ws.Range("N34").Value = -5005.4

# Row 92
$ws.Range("H92").Value = 49239.2
$ws.Range("J92").Value = 49239.2
$ws.Range("L92").Value = 49239.2
$ws.Range("N92").Value = -54231.2

# Row 122
$ws.Range("H122").Value = 2000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 2682.0908
$ws.Range("I132").Value = 2253.2222
$ws.Range("K132").Value = 6759.6666
$ws.Range("M132").Value = -4229.6666

# Row 134
$ws.Range("H134").Value = 800
$ws.Range("I134").Value = 800
$ws.Range("K134").Value = 2400
$ws.Range("M134").Value = 135


# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 107
$ws.Range("H107").Value = 2177.5
$ws.Range("I107").Value = 1300
$ws.Range("J107").Value = 3932.5
$ws.Range("K107").Value = 1300
$ws.Range("L107").Value = 3932.5
$ws.Range("M107").Value = 620
$ws.Range("N107").Value = -7772.5

# Row 132
$ws.Range("H132").Value = 2916.75
$ws.Range("I132").Value = 2444.7778
$ws.Range("K132").Value = 7334.3334
$ws.Range("M132").Value = -4804.3334


# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 3622.75
$ws.Range("I46").Value = 3622.75
$ws.Range("K46").Value = 3622.75
$ws.Range("M46").Value = -3434.75

# Row 55
$ws.Range("H55").Value = 11111
$ws.Range("I55").Value = 11111
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 11111
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = -10938
$ws.Range("M55").ClearContents()

# Row 61
$ws.Range("H61").Value = 1098.75
$ws.Range("I61").Value = 1098.75
$ws.Range("K61").Value = 1098.75
$ws.Range("M61").Value = -896.75

# Row 113
$ws.Range("H113").Value = 1098.75
$ws.Range("I113").Value = 1098.75
$ws.Range("K113").Value = 1098.75
$ws.Range("M113").Value = 1071.25

# Row 132
$ws.Range("H132").Value = 2689.7273
$ws.Range("I132").Value = 1954.3334
$ws.Range("K132").Value = 5863.0002
$ws.Range("M132").Value = -3333.0002


# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 38
$ws.Range("H38").Value = 28990
$ws.Range("J38").Value = 28990
$ws.Range("L38").Value = 28990
$ws.Range("N38").Value = -29936

# Row 113
$ws.Range("H113").Value = 855.5
$ws.Range("I113").Value = 878.4286
$ws.Range("K113").Value = 2635.2858
$ws.Range("M113").Value = -465.2857999999997

# Row 122
$ws.Range("H122").Value = 2320.875
$ws.Range("J122").Value = 2193.75
$ws.Range("L122").Value = 6581.25
$ws.Range("N122").Value = -11481.25

# Row 132
$ws.Range("H132").Value = 4054.7036
$ws.Range("I132").Value = 3292.25
$ws.Range("J132").Value = 5163.727
$ws.Range("K132").Value = 9876.75
$ws.Range("L132").Value = 15491.181
$ws.Range("M132").Value = -7346.75
$ws.Range("N132").Value = -20551.181

# Row 136
$ws.Range("H136").Value = 901.7143
$ws.Range("I136").Value = 901.5
$ws.Range("K136").Value = 2704.5
$ws.Range("M136").Value = -154.5

